$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'kimchi mix',
    'kinesiology device',
    'kinesiology gifts',
    'kinesiology knee brace',
    'kinesiology rock tape',
    'kinesiology shirt',
    'kinesiology spartan tape',
    'kinesiology sports tape',
    'kinesiology tape 4',
    'kinesiology tape 4 inch',
    'kinesiology tape 4 pack',
    'kinesiology tape achilles',
    'kinesiology tape for sensitive skin',
    'kinesiology tape gentle',
    'kinesiology tape guide',
    'kinesiology tape hamstring',
    'kinesiology tape hip',
    'kinesiology tape knee',
    'kinesiology tape lower back',
    'kinesiology tape ok tape',
    'kinesiology tape pro',
    'kinesiology tape red',
    'kinesiology tape set',
    'kinesiology tape spartan',
    'kinesiology tape x',
    'kingsland equestrian apparel women',
    'knee air brace',
    'knee and hip joint pain',
    'knee and joint',
    'knee band basketball',
    'knee band for running',
    'knee band running',
    'knee band support',
    'knee bands for runners knee',
    'knee bands for running',
    'knee bands running',
    'knee black leggings',
    'knee bleach',
    'knee brace 3x',
    'knee brace 5x',
    'knee brace for athletes',
    'knee brace for larger legs',
    'knee brace for runners knee',
    'knee brace for running pair',
    'knee brace for running women',
    'knee brace for running xl',
    'knee brace for skiing',
    'knee brace for standing all day',
    'knee brace hiking',
    'knee brace hot cold',
    'knee brace hyperextension',
    'knee brace leggings',
    'knee brace marathon',
    'knee brace over pants',
    'knee brace patella tracking',
    'knee brace plus',
    'knee brace size 2x',
    'knee brace size chart',
    'knee brace skiing',
    'knee brace skin color',
    'knee brace stability',
    'knee brace tape',
    'knee braces for runners knee women',
    'knee braces for skiing',
    'knee compression',
    'knee compression basketball',
    'knee compression cold',
    'knee compression crossfit',
    'knee compression for running',
    'knee compression for women',
    'knee compression garment',
    'knee compression leggings',
    'knee compression machine',
    'knee compression pants',
    'knee compression plus size',
    'knee compression running',
    'knee compression shorts',
    'knee compression support',
    'knee compression tape',
    'knee compression tights',
    'knee compression women',
    'knee compression wraps female',
    'knee compressions',
    'knee dislocation',
    'knee for running',
    'knee gym support',
    'knee high compression',
    'knee high for women',
    'knee high leggings',
    'knee high leggings for women',
    'knee high pants for women',
    'knee high panty hose',
    'knee high panty hose black',
    'knee high pany hose',
    'knee high shorts',
    'knee high shorts for women',
    'knee high stocking',
    'knee high stockings black',
    'knee high support',
    'knee high tights'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

